# "new tasks + switch"
# - Row 3 "Draw Bar" toolbar: rename/relabel a couple of buttons and recolor
#   the button strip (B3/F3 "switch <=>" -> "switch ", D3 becomes "draw clr",
#   B3 becomes the freed-up "switch " label) plus a small border tidy-up.
# - Add a brand-new "Tool Bar" section: a merged title row (A4:H4, styled
#   like the existing "Draw Bar" title in A2:H2) followed by a populated
#   button strip in row 5 (copy/cut/select figure/save/paste/delete/
#   save (type)/load), styled like the row-3 strip.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Global font refresh (workbook-wide font swap captured in the diff)
# ---------------------------------------------------------------------
$ws.Cells.Font.Name = "Calibri"

# ---------------------------------------------------------------------
# Row 3 toolbar: relabel + recolor the affected buttons
# ---------------------------------------------------------------------
# B3 was "draw clr" (yellow fill clr-swatch) -> becomes "switch " (teal, like
# the rest of the strip, thick right divider)
$ws.Range("B3").Value = "switch "
$ws.Range("B3").Interior.ThemeColor = 9
$ws.Range("B3").Borders.Item(10).LineStyle = 1
$ws.Range("B3").Borders.Item(10).Weight = 4

# D3 was "select figure" -> becomes "draw clr" (red/accent2 fill clr-swatch,
# thick right divider)
$ws.Range("D3").Value = "draw clr"
$ws.Range("D3").Interior.ThemeColor = 6
$ws.Range("D3").Borders.Item(10).LineStyle = 1
$ws.Range("D3").Borders.Item(10).Weight = 4

# F3 label simplified from "switch <=>" to "switch " ; thick right divider
$ws.Range("F3").Value = "switch "
$ws.Range("F3").Borders.Item(10).LineStyle = 1
$ws.Range("F3").Borders.Item(10).Weight = 4

# H3 "fill clr" swatch switches from yellow to red/accent2, right divider
# goes back to a plain thin line (it's the last cell before the legend box)
$ws.Range("H3").Interior.ThemeColor = 6
$ws.Range("H3").Borders.Item(10).LineStyle = 1
$ws.Range("H3").Borders.Item(10).Weight = 2

# Minor border tidy-up: I3 loses its stray left border
$ws.Range("I3").Borders.Item(7).LineStyle = 0

# ---------------------------------------------------------------------
# New "Tool Bar" title (row 4), mirrors the "Draw Bar" title in row 2
# ---------------------------------------------------------------------
$title = $ws.Range("A4:H4")
$title.Interior.ThemeColor = 4
$title.Borders.Item(8).LineStyle = 1
$title.Borders.Item(8).ThemeColor = 2
$title.Font.Size = 22
$title.Font.Name = "Calibri"
$title.HorizontalAlignment = -4108
$title.VerticalAlignment = -4108
$title.Merge()
$ws.Range("A4").Value = "Tool Bar"
$ws.Rows.Item(4).RowHeight = 28.5

# ---------------------------------------------------------------------
# New Tool Bar button strip (row 5), same template as row 3
# ---------------------------------------------------------------------
$row5Values = @{
    "A" = "copy"
    "B" = "cut"
    "C" = "select figure"
    "D" = "save"
    "E" = "paste"
    "F" = "delete"
    "G" = "save (type)"
    "H" = "load"
}
# C column keeps the teal (accent5) swatch, everything else is red (accent2)
$tealCols = @("C")
$allCols = @("A","B","C","D","E","F","G","H")
$boxCols = @("A","H")
$thickRightCols = @("B","D","F")

foreach ($col in $allCols) {
    $cell = $ws.Range($col + "5")
    $cell.Value = $row5Values[$col]
    $cell.Font.Name = "Calibri"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108

    if ($tealCols -contains $col) {
        $cell.Interior.ThemeColor = 9
    } else {
        $cell.Interior.ThemeColor = 6
    }

    # top + bottom thin border on every cell
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1

    if ($boxCols -contains $col) {
        # full thin box
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(10).LineStyle = 1
    } elseif ($thickRightCols -contains $col) {
        # thin left, thick right divider
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight = 4
    } else {
        # no left border, thin right border
        $cell.Borders.Item(10).LineStyle = 1
    }
}

# ---------------------------------------------------------------------
# Restore the selection to where the author left off
# ---------------------------------------------------------------------
$ws.Range("H3").Select()
